$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "244.89"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.88"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.392"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8151"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9351"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1436"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07446"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03522"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03068"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09417"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.012"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04797"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0005940"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.005587"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.004167"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0009875"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.670"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.416"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1347"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.00006998"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002900"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04004"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006348"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1077"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002899"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.005918"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005238"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.0000"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002367"
